$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp in the title row ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 10:35"

# --- Helper: locate a country's data row by name and update its stats ---
function Set-CountryStats {
    param($Country, $TotalCases, $NewCases, $ActiveCases, $Recovered, $Critical, $DeathsToday, $Deaths)

    $found = $ws.Range("A4:A219").Find($Country)
    $r = $found.Row

    $ws.Cells.Item($r, 2).Value = $TotalCases
    $ws.Cells.Item($r, 3).Value = $NewCases
    $ws.Cells.Item($r, 4).Value = $ActiveCases
    $ws.Cells.Item($r, 5).Value = $Recovered
    $ws.Cells.Item($r, 6).Value = $Critical
    $ws.Cells.Item($r, 7).Value = $DeathsToday
    $ws.Cells.Item($r, 8).Value = $Deaths
}

# Polonia: new totals push it above Indonesia
Set-CountryStats "Polonia" 20838 219 8977 10871 0 8 990

# Filipinas: updated totals (stays in the same sorted position)
Set-CountryStats "Filipinas" 13777 180 3177 9737 0 6 863

# Afganistan: updated totals (stays in the same sorted position)
Set-CountryStats "Afganistan" 9998 782 1040 8742 0 11 216

# Namibia: new totals push it above Laos
Set-CountryStats "Namibia" 20 1 14 6 0 0 0

# Sri Lanka: only active cases / recovered change
$found = $ws.Range("A4:A219").Find("Sri Lanka")
$r = $found.Row
$ws.Cells.Item($r, 4).Value = 660
$ws.Cells.Item($r, 5).Value = 399

# --- Re-sort the country table by "Casos totales" (column B) descending,
#     matching the workbook's existing ordering convention ---
$dataRange = $ws.Range("A4:H219")
$dataRange.Sort($ws.Range("B4"), 2)
